$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-10-08 12:16:31", "check_availability", "https://example.com", "Checked availability: Selected or default date current date is available for booking.", "2024-10-08", "12:16:31"),
    @("2024-10-08 12:16:32", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-10-08", "12:16:32"),
    @("2024-10-08 12:16:32", "check_availability", "https://example.com", "Checked availability: No availability for the selected date.", "2024-10-08", "12:16:32"),
    @("2024-10-08 12:16:33", "check_availability", "https://example.com", "Checked availability: Selected or default date is available for booking.", "2024-10-08", "12:16:33"),
    @("2024-10-08 12:16:34", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-10-08", "12:16:34")
)

$startRow = 326
$endRow = $startRow + $rows.Count - 1
$writeRange = $ws.Range("A$startRow`:F$endRow")
$writeRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}

$writeRange.ClearFormats()
